# Applies the Q4-2018/Q3-2018 financial-data refresh described in the commit:
# inserts two new quarter columns (D:E) ahead of the existing data and restates
# a handful of older quarters whose figures changed in the source refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two blank columns before column D (old D:K shifts right to F:M)
$ws.Range("D5:E102").Insert(-4161)

# New columns should look like the old D/E (number/date formatting) columns,
# which have shifted right to F/G - copy their formats over.
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Period-ending dates for the two newly inserted quarter columns (Dec-2018, Sep-2018)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373

# Populate the new-quarter values for every data row
$ws.Range("D8").Value = 404300
$ws.Range("E8").Value = 429900
$ws.Range("D9").Value = 367400
$ws.Range("E9").Value = 387000
$ws.Range("D10").Value = 36900
$ws.Range("E10").Value = 42900
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 391600
$ws.Range("E17").Value = 409000
$ws.Range("D18").Value = 12700
$ws.Range("E18").Value = 20900
$ws.Range("D20").Value = 400
$ws.Range("E20").Value = 200
$ws.Range("D21").Value = 18700
$ws.Range("E21").Value = 26500
$ws.Range("D22").Value = 4300
$ws.Range("E22").Value = 4300
$ws.Range("D23").Value = 8800
$ws.Range("E23").Value = 16800
$ws.Range("D24").Value = 2100
$ws.Range("E24").Value = 1900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 6700
$ws.Range("E26").Value = 14900
$ws.Range("D27").Value = 6600
$ws.Range("E27").Value = 14800
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -400
$ws.Range("E32").Value = -200
$ws.Range("D33").Value = 6600
$ws.Range("E33").Value = 14800
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 6600
$ws.Range("E35").Value = 14800
$ws.Range("D41").Value = 125500
$ws.Range("E41").Value = 101100
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 168400
$ws.Range("E43").Value = 203600
$ws.Range("D44").Value = 218200
$ws.Range("E44").Value = 218800
$ws.Range("D45").Value = 8500
$ws.Range("E45").Value = 10200
$ws.Range("D46").Value = 520600
$ws.Range("E46").Value = 533700
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 147800
$ws.Range("E48").Value = 143400
$ws.Range("D49").Value = 6000
$ws.Range("E49").Value = 6100
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 16600
$ws.Range("E52").Value = 15400
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 691000
$ws.Range("E54").Value = 698600
$ws.Range("D57").Value = 114100
$ws.Range("E57").Value = 121200
$ws.Range("D58").Value = 4600
$ws.Range("E58").Value = 5000
$ws.Range("D59").Value = 40300
$ws.Range("E59").Value = 31700
$ws.Range("D60").Value = 159000
$ws.Range("E60").Value = 157900
$ws.Range("D61").Value = 305700
$ws.Range("E61").Value = 306300
$ws.Range("D62").Value = 38500
$ws.Range("E62").Value = 38200
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 508200
$ws.Range("E66").Value = 507200
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 148800
$ws.Range("E72").Value = 144300
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 182800
$ws.Range("E76").Value = 191400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D81").Value = 6600
$ws.Range("E81").Value = 14800
$ws.Range("D83").Value = 5600
$ws.Range("E83").Value = 5400
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 48800
$ws.Range("E89").Value = 18000
$ws.Range("D91").Value = -8800
$ws.Range("E91").Value = -4700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -8700
$ws.Range("E94").Value = -4800
$ws.Range("D96").Value = -2000
$ws.Range("E96").Value = -1900
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -15800
$ws.Range("E100").Value = -2600
$ws.Range("D101").Value = 100
$ws.Range("E101").Value = -600
$ws.Range("D102").Value = 24400
$ws.Range("E102").Value = 10000

# Restatements to a handful of older quarters (figures revised in this refresh)
$ws.Range("H8").Value = 424900
$ws.Range("I8").Value = 359400
$ws.Range("H9").Value = 379100
$ws.Range("I9").Value = 315800
$ws.Range("H10").Value = 45800
$ws.Range("I10").Value = 43600
$ws.Range("G14").Value = "NA"
$ws.Range("G17").Value = 446200
$ws.Range("H17").Value = 401300
$ws.Range("I17").Value = 335400
$ws.Range("G18").Value = 25600
$ws.Range("H18").Value = 23600
$ws.Range("I18").Value = 24000
$ws.Range("G20").Value = -200
$ws.Range("H21").Value = 28300
$ws.Range("I21").Value = 27800
$ws.Range("H23").Value = 19300
$ws.Range("I23").Value = 19000
$ws.Range("I24").Value = 6600
$ws.Range("H26").Value = 13100
$ws.Range("I26").Value = 12400
$ws.Range("H27").Value = 12900
$ws.Range("I27").Value = 12300
$ws.Range("G32").Value = 200
$ws.Range("H33").Value = 5600
$ws.Range("I33").Value = 12300
$ws.Range("H35").Value = 5600
$ws.Range("I35").Value = 12300
$ws.Range("H81").Value = 5600
$ws.Range("I81").Value = 12300

# Workbook-level calculation option tweak from the refresh tool
$wb.Application.Calculation = -4135
